$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 8418
$ws.Range("F7").Value = 143
$ws.Range("F11").Value = 7687
$ws.Range("F12").Value = 7862
$ws.Range("F13").Value = 5041
$ws.Range("F17").Value = 5407
$ws.Range("F19").Value = 620
$ws.Range("F20").Value = 147
$ws.Range("F22").Value = 974
$ws.Range("F23").Value = 1511
$ws.Range("F24").Value = 2121
$ws.Range("F27").Value = 265
$ws.Range("F32").Value = 802
$ws.Range("F33").Value = 1308
$ws.Range("F34").Value = 487

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F31").Value = 84

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 1648
$ws.Range("F7").Value = 698
$ws.Range("F9").Value = 9495
$ws.Range("F10").Value = 1824
$ws.Range("G12").Value = "已售罄"
$ws.Range("F15").Value = 312
$ws.Range("F16").Value = 2622
$ws.Range("F17").Value = 298
$ws.Range("F18").Value = 112

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 1824
$ws.Range("F7").Value = 312
$ws.Range("F8").Value = 2622
$ws.Range("F9").Value = 298
$ws.Range("F11").Value = 7687
$ws.Range("F12").Value = 7862
$ws.Range("F14").Value = 620
$ws.Range("F15").Value = 147
$ws.Range("F16").Value = 112
$ws.Range("F17").Value = 974
$ws.Range("F18").Value = 1511
$ws.Range("F19").Value = 2121
$ws.Range("F26").Value = 265
$ws.Range("F29").Value = 802
$ws.Range("F31").Value = 1308
$ws.Range("F35").Value = 487
